# Liêm Trinh Hoàn Thiện
# Adds new "Liêm Trinh" (Mệnh cung) lookup rows to the table in columns A:B,
# rows 130-173, matching style/pattern of the existing rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
  ,@(130, "Liêm Trinh tọa thủ cung Mệnh ở Tỵ gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình", "Liêm Trinh tọa thủ cung Mệnh ở Tỵ gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình")
  ,@(131, "Liêm Trinh tọa thủ cung Mệnh ở Hợi gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình", "Liêm Trinh tọa thủ cung Mệnh ở Hợi gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình")
  ,@(132, "Liêm Trinh tọa thủ cung Mệnh ở Tỵ gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình", "Liêm Trinh tọa thủ cung Mệnh ở Tỵ gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình")
  ,@(133, "Liêm Trinh tọa thủ cung Mệnh ở Hợi gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình", "Liêm Trinh tọa thủ cung Mệnh ở Hợi gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình")
  ,@(134, "Liêm Trinh tọa thủ cung Mệnh ở Mão gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình", "Liêm Trinh tọa thủ cung Mệnh ở Mão gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình")
  ,@(135, "Liêm Trinh tọa thủ cung Mệnh ở Dậu gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình", "Liêm Trinh tọa thủ cung Mệnh ở Dậu gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình")
  ,@(136, "Liêm Trinh tọa thủ cung Mệnh ở Tỵ", "Liêm Trinh tọa thủ cung Mệnh ở Tỵ")
  ,@(137, "Liêm Trinh tọa thủ cung Mệnh ở Hợi gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình", "Liêm Trinh tọa thủ cung Mệnh ở Hợi")
  ,@(138, "Liêm Trinh tọa thủ cung Mệnh ở Mão gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình", "Liêm Trinh tọa thủ cung Mệnh ở Mão")
  ,@(139, "Liêm Trinh tọa thủ cung Mệnh ở Hợi", "Liêm Trinh tọa thủ cung Mệnh ở Hợi")
  ,@(140, "Liêm Trinh đồng cung Hoá Kỵ tại Tỵ, Văn Xương, Văn Khúc tại Mệnh và sinh năm Bính", "Liêm Trinh đồng cung Hoá Kỵ tại Tỵ, Văn Xương, Văn Khúc tại Mệnh và sinh năm Bính")
  ,@(141, "Liêm Trinh đồng cung Hoá Kỵ tại Hợi, Văn Xương, Văn Khúc tại Mệnh và sinh năm Bính", "Liêm Trinh đồng cung Hoá Kỵ tại Hợi, Văn Xương, Văn Khúc tại Mệnh và sinh năm Bính")
  ,@(142, "Liêm Trinh tọa thủ cung Mệnh ở Dậu gặp Hoả Tinh, Linh Tinh", "Liêm Trinh tọa thủ cung Mệnh ở Dậu gặp Hoả Tinh, Linh Tinh")
  ,@(143, "Liêm Trinh tọa thủ cung Mệnh ở Mão gặp Hoả Tinh, Linh Tinh", "Liêm Trinh tọa thủ cung Mệnh ở Mão gặp Hoả Tinh, Linh Tinh")
  ,@(144, "Anh có Liêm Trinh tọa thủ cung Mệnh ở Tý gặp Văn Xương, Văn Khúc", "Anh có Liêm Trinh tọa thủ cung Mệnh ở Tý gặp Văn Xương, Văn Khúc")
  ,@(145, "Anh có Liêm Trinh tọa thủ cung Mệnh ở Ngọ gặp Văn Xương, Văn Khúc", "Anh có Liêm Trinh tọa thủ cung Mệnh ở Ngọ gặp Văn Xương, Văn Khúc")
  ,@(146, "Anh có Liêm Trinh tọa thủ cung Mệnh ở Dần gặp Văn Xương, Văn Khúc", "Anh có Liêm Trinh tọa thủ cung Mệnh ở Dần gặp Văn Xương, Văn Khúc")
  ,@(147, "Anh có Liêm Trinh tọa thủ cung Mệnh ở Thân gặp Văn Xương, Văn Khúc", "Anh có Liêm Trinh tọa thủ cung Mệnh ở Thân gặp Văn Xương, Văn Khúc")
  ,@(148, "Anh có Liêm Trinh tọa thủ cung Mệnh ở Thìn gặp Văn Xương, Văn Khúc", "Anh có Liêm Trinh tọa thủ cung Mệnh ở Thìn gặp Văn Xương, Văn Khúc")
  ,@(149, "Anh có Liêm Trinh tọa thủ cung Mệnh ở Tuất gặp Văn Xương, Văn Khúc", "Anh có Liêm Trinh tọa thủ cung Mệnh ở Tuất gặp Văn Xương, Văn Khúc")
  ,@(150, "Anh có Liêm Trinh tọa thủ cung Mệnh ở Sửu gặp Văn Xương, Văn Khúc", "Anh có Liêm Trinh tọa thủ cung Mệnh ở Sửu gặp Văn Xương, Văn Khúc")
  ,@(151, "Anh có Liêm Trinh tọa thủ cung Mệnh ở Mùi gặp Văn Xương, Văn Khúc", "Anh có Liêm Trinh tọa thủ cung Mệnh ở Mùi gặp Văn Xương, Văn Khúc")
  ,@(152, "Chị có Liêm Trinh tọa thủ cung Mệnh ở Tý gặp Văn Xương, Văn Khúc", "Chị có Liêm Trinh tọa thủ cung Mệnh ở Tý gặp Văn Xương, Văn Khúc")
  ,@(153, "Chị có Liêm Trinh tọa thủ cung Mệnh ở Ngọ gặp Văn Xương, Văn Khúc", "Chị có Liêm Trinh tọa thủ cung Mệnh ở Ngọ gặp Văn Xương, Văn Khúc")
  ,@(154, "Chị có Liêm Trinh tọa thủ cung Mệnh ở Dần gặp Văn Xương, Văn Khúc", "Chị có Liêm Trinh tọa thủ cung Mệnh ở Dần gặp Văn Xương, Văn Khúc")
  ,@(155, "Chị có Liêm Trinh tọa thủ cung Mệnh ở Thân gặp Văn Xương, Văn Khúc", "Chị có Liêm Trinh tọa thủ cung Mệnh ở Thân gặp Văn Xương, Văn Khúc")
  ,@(156, "Chị có Liêm Trinh tọa thủ cung Mệnh ở Thìn gặp Văn Xương, Văn Khúc", "Chị có Liêm Trinh tọa thủ cung Mệnh ở Thìn gặp Văn Xương, Văn Khúc")
  ,@(157, "Chị có Liêm Trinh tọa thủ cung Mệnh ở Tuất gặp Văn Xương, Văn Khúc", "Chị có Liêm Trinh tọa thủ cung Mệnh ở Tuất gặp Văn Xương, Văn Khúc")
  ,@(158, "Chị có Liêm Trinh tọa thủ cung Mệnh ở Sửu gặp Văn Xương, Văn Khúc", "Chị có Liêm Trinh tọa thủ cung Mệnh ở Sửu gặp Văn Xương, Văn Khúc")
  ,@(159, "Chị có Liêm Trinh tọa thủ cung Mệnh ở Mùi gặp Văn Xương, Văn Khúc", "Chị có Liêm Trinh tọa thủ cung Mệnh ở Mùi gặp Văn Xương, Văn Khúc")
  ,@(160, "Liêm Trinh đồng cung Bạch Hổ tại Tỵ", "Liêm Trinh đồng cung Bạch Hổ tại Tỵ")
  ,@(161, "Liêm Trinh đồng cung Bạch Hổ tại Hợi", "Liêm Trinh đồng cung Bạch Hổ tại Hợi")
  ,@(162, "Liêm Trinh đồng cung Bạch Hổ tại Mão", "Liêm Trinh đồng cung Bạch Hổ tại Mão")
  ,@(163, "Liêm Trinh đồng cung Bạch Hổ tại Dậu", "Liêm Trinh đồng cung Bạch Hổ tại Dậu")
  ,@(164, "Liêm Trinh đồng cung Thiên Tướng tại Thân", "Liêm Trinh đồng cung Thiên Tướng tại Thân")
  ,@(165, "Liêm Trinh đồng cung Thiên Tướng tại Dần", "Liêm Trinh đồng cung Thiên Tướng tại Dần")
  ,@(166, "Liêm Trinh đồng cung Thiên Tướng tại Ngọ", "Liêm Trinh đồng cung Thiên Tướng tại Ngọ")
  ,@(167, "Liêm Trinh đồng cung Thiên Tướng tại Tý", "Liêm Trinh đồng cung Thiên Tướng tại Tý")
  ,@(168, "Liêm Trinh đồng cung Thiên Tướng tại Mùi", "Liêm Trinh đồng cung Thiên Tướng tại Mùi")
  ,@(169, "Liêm Trinh đồng cung Thiên Tướng tại Sửu", "Liêm Trinh đồng cung Thiên Tướng tại Sửu")
  ,@(170, "Liêm Trinh đồng cung Thiên Tướng tại Tuất", "Liêm Trinh đồng cung Thiên Tướng tại Tuất")
  ,@(171, "Liêm Trinh đồng cung Thiên Tướng tại Thìn", "Liêm Trinh đồng cung Thiên Tướng tại Thìn")
  ,@(172, "Liêm Trinh tọa thủ cung Mệnh ở Mão gặp Hoả Tinh, Linh Tinh", "Liêm Trinh tọa thủ cung Mệnh ở Mão gặp Hoả Tinh, Linh Tinh")
  ,@(173, "Liêm Trinh tọa thủ cung Mệnh ở Dậu gặp Hoả Tinh, Linh Tinh", "Liêm Trinh tọa thủ cung Mệnh ở Dậu gặp Hoả Tinh, Linh Tinh")
)

foreach ($item in $rowsData) {
  $r = $item[0]
  $a = $item[1]
  $b = $item[2]
  $ws.Cells.Item($r, 1).Value2 = $a
  $ws.Cells.Item($r, 2).Value2 = $b
  $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 2)).Interior.Color = 65535
}

# Restore the view state captured in the saved workbook: scrolled down to
# row 142 with K158 as the active (selected) cell.
$ws.Activate()
try {
  $excel.ActiveWindow.ScrollRow = 142
  $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("K158").Select()

Write-Output "Added rows 130-173 (Liêm Trinh entries)."
